$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 393.625
$ws.Range("I2").Value = 338.85715
$ws.Range("K2").Value = 338.85715
$ws.Range("M2").Value = -225.85715
$ws.Range("H3").Value = 43657
$ws.Range("J3").Value = 43657
$ws.Range("L3").Value = 43657
$ws.Range("N3").Value = -43885
$ws.Range("H6").Value = 33336092
$ws.Range("J6").Value = 4998.875
$ws.Range("L6").Value = 14996.625
$ws.Range("N6").Value = -15220.625
$ws.Range("H17").Value = 960.14
$ws.Range("J17").Value = 1029.0222
$ws.Range("L17").Value = 3087.0666
$ws.Range("N17").Value = -3423.0666
$ws.Range("H20").Value = 906.6667
$ws.Range("I20").Value = 906.6667
$ws.Range("K20").Value = 906.6667
$ws.Range("M20").Value = -676.6667
$ws.Range("H29").Value = 71428696
$ws.Range("I29").Value = 71428696
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 214286088
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -214285807
$ws.Range("N29").ClearContents()
$ws.Range("H34").Value = 3000
$ws.Range("I34").Value = 3000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2797
$ws.Range("H35").Value = 906.6667
$ws.Range("I35").Value = 906.6667
$ws.Range("K35").Value = 906.6667
$ws.Range("M35").Value = -527.6667
$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 3000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2285
$ws.Range("H38").Value = 11428.143
$ws.Range("I38").Value = 16285.429
$ws.Range("K38").Value = 48856.287
$ws.Range("M38").Value = -48484.287
$ws.Range("H39").Value = 1797.4
$ws.Range("I39").Value = 999.5
$ws.Range("K39").Value = 2998.5
$ws.Range("M39").Value = -2702.5
$ws.Range("H40").Value = 3017.5715
$ws.Range("I40").Value = 1863.6364
$ws.Range("J40").Value = 3764.2354
$ws.Range("K40").Value = 1863.6364
$ws.Range("L40").Value = 3764.2354
$ws.Range("M40").Value = -1688.6364
$ws.Range("N40").Value = -4114.2354
$ws.Range("H43").Value = 5838
$ws.Range("I43").Value = 1416
$ws.Range("J43").Value = 7496.25
$ws.Range("K43").Value = 1416
$ws.Range("L43").Value = 7496.25
$ws.Range("M43").Value = -1347
$ws.Range("N43").Value = -7634.25
$ws.Range("H51").Value = 3086.1738
$ws.Range("J51").Value = 3334.1765
$ws.Range("L51").Value = 3334.1765
$ws.Range("N51").Value = -4302.1765
$ws.Range("H53").Value = 490.70834
$ws.Range("I53").Value = 379.66666
$ws.Range("J53").Value = 675.7778
$ws.Range("K53").Value = 379.66666
$ws.Range("L53").Value = 675.7778
$ws.Range("M53").Value = 257.33334
$ws.Range("N53").Value = -1949.7778
$ws.Range("H58").Value = 62500856
$ws.Range("I58").Value = 83333480
$ws.Range("K58").Value = 250000440
$ws.Range("M58").Value = -250000290
$ws.Range("H74").Value = 3157.5557
$ws.Range("I74").Value = 3157.5557
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3157.5557
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2221.5557
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 3147.25
$ws.Range("I76").Value = 3189
$ws.Range("K76").Value = 3189
$ws.Range("M76").Value = -2874
$ws.Range("H77").Value = 3157.5557
$ws.Range("I77").Value = 3157.5557
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 15787.7785
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -11107.7785
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 3147.25
$ws.Range("I79").Value = 3189
$ws.Range("K79").Value = 3189
$ws.Range("M79").Value = -2097
$ws.Range("H86").Value = 4271.222
$ws.Range("I86").Value = 4500
$ws.Range("J86").Value = 4088.2
$ws.Range("K86").Value = 4500
$ws.Range("L86").Value = 4088.2
$ws.Range("M86").Value = -3377
$ws.Range("N86").Value = -6334.2
$ws.Range("H89").Value = 4271.222
$ws.Range("I89").Value = 4500
$ws.Range("J89").Value = 4088.2
$ws.Range("K89").Value = 22500
$ws.Range("L89").Value = 20441
$ws.Range("M89").Value = -16884
$ws.Range("N89").Value = -31673
$ws.Range("H98").Value = 3052.158
$ws.Range("I98").Value = 3008.5454
$ws.Range("J98").Value = 3112.125
$ws.Range("K98").Value = 3008.5454
$ws.Range("L98").Value = 3112.125
$ws.Range("M98").Value = -1510.5454
$ws.Range("N98").Value = -6108.125
$ws.Range("H102").Value = 43657
$ws.Range("J102").Value = 43657
$ws.Range("L102").Value = 43657
$ws.Range("N102").Value = -50147
$ws.Range("H106").Value = 5424.6665
$ws.Range("I106").Value = 5117.5713
$ws.Range("J106").Value = 6499.5
$ws.Range("K106").Value = 5117.5713
$ws.Range("L106").Value = 6499.5
$ws.Range("M106").Value = -4486.5713
$ws.Range("N106").Value = -7761.5
$ws.Range("H116").Value = 4205.625
$ws.Range("J116").Value = 4502.75
$ws.Range("L116").Value = 4502.75
$ws.Range("N116").Value = -11386.75
$ws.Range("H122").Value = 3052.158
$ws.Range("I122").Value = 3008.5454
$ws.Range("J122").Value = 3112.125
$ws.Range("K122").Value = 9025.636200000001
$ws.Range("L122").Value = 9336.375
$ws.Range("M122").Value = -6575.636200000001
$ws.Range("N122").Value = -14236.375
$ws.Range("H124").Value = 50709
$ws.Range("I124").Value = 50709
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 50709
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -45799
$ws.Range("N124").ClearContents()
$ws.Range("H129").Value = 2076.6
$ws.Range("J129").Value = 3239.2
$ws.Range("L129").Value = 9717.599999999999
$ws.Range("N129").Value = -19717.6
$ws.Range("H132").Value = 2059.7576
$ws.Range("J132").Value = 3524.818
$ws.Range("L132").Value = 10574.454
$ws.Range("N132").Value = -15634.454
$ws.Range("H135").Value = 1781.1111
$ws.Range("I135").Value = 1781.1111
$ws.Range("K135").Value = 16029.9999
$ws.Range("M135").Value = -13494.9999
$ws.Range("H137").Value = 6280.65
$ws.Range("I137").Value = 6977.0586
$ws.Range("K137").Value = 20931.1758
$ws.Range("M137").Value = -18381.1758
$ws.Range("H138").Value = 5850.344
$ws.Range("I138").Value = 1174.6666
$ws.Range("J138").Value = 6569.6797
$ws.Range("K138").Value = 3523.9998
$ws.Range("L138").Value = 19709.0391
$ws.Range("M138").Value = 1616.0002
$ws.Range("N138").Value = -29989.0391
$ws.Range("H141").Value = 8122.3
$ws.Range("I141").Value = 6392.5
$ws.Range("J141").Value = 8863.643
$ws.Range("K141").Value = 19177.5
$ws.Range("L141").Value = 26590.929
$ws.Range("M141").Value = -13997.5
$ws.Range("N141").Value = -36950.929

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 6089.6665
$ws.Range("I61").Value = 6089.6665
$ws.Range("K61").Value = 6089.6665
$ws.Range("M61").Value = -5877.6665
$ws.Range("H88").Value = 3510.7
$ws.Range("I88").Value = 2303
$ws.Range("K88").Value = 2303
$ws.Range("M88").Value = -1897
$ws.Range("H91").Value = 3510.7
$ws.Range("I91").Value = 2303
$ws.Range("K91").Value = 2303
$ws.Range("M91").Value = -899
$ws.Range("H132").Value = 2146.4443
$ws.Range("I132").Value = 2059.7917
$ws.Range("J132").Value = 2839.6667
$ws.Range("K132").Value = 6179.375100000001
$ws.Range("L132").Value = 8519.000100000001
$ws.Range("M132").Value = -3649.375100000001
$ws.Range("N132").Value = -13579.0001
$ws.Range("H136").Value = 6089.6665
$ws.Range("I136").Value = 6089.6665
$ws.Range("K136").Value = 18268.9995
$ws.Range("M136").Value = -15718.9995

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H19").Value = 3875
$ws.Range("H20").Value = 2179.25
$ws.Range("J20").Value = 1521.25
$ws.Range("L20").Value = 1521.25
$ws.Range("N20").Value = -2015.25
$ws.Range("H75").Value = 19995
$ws.Range("I75").Value = 19995
$ws.Range("K75").Value = 19995
$ws.Range("M75").Value = -19059
$ws.Range("H78").Value = 19995
$ws.Range("I78").Value = 19995
$ws.Range("K78").Value = 59985
$ws.Range("M78").Value = -55305
$ws.Range("H86").Value = 25002628
$ws.Range("I86").Value = 2473.6875
$ws.Range("K86").Value = 2473.6875
$ws.Range("M86").Value = -1350.6875
$ws.Range("H89").Value = 25002628
$ws.Range("I89").Value = 2473.6875
$ws.Range("K89").Value = 12368.4375
$ws.Range("M89").Value = -6752.4375
$ws.Range("H94").Value = 1801.0817
$ws.Range("I94").Value = 1609.2778
$ws.Range("K94").Value = 1609.2778
$ws.Range("M94").Value = -1158.2778

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 4089.25
$ws.Range("I31").Value = 1591
$ws.Range("K31").Value = 1591
$ws.Range("M31").Value = -1296
$ws.Range("H34").Value = 4089.25
$ws.Range("I34").Value = 1591
$ws.Range("K34").Value = 1591
$ws.Range("M34").Value = -1389
$ws.Range("H41").Value = 19803.666
$ws.Range("J41").Value = 53659.668
$ws.Range("L41").Value = 53659.668
$ws.Range("N41").Value = -54515.668
$ws.Range("H50").Value = 49744.5
$ws.Range("J50").Value = 49744.5
$ws.Range("L50").Value = 49744.5
$ws.Range("N50").Value = -50994.5
$ws.Range("H51").Value = 19099.092
$ws.Range("H56").Value = 28966.334
$ws.Range("J56").Value = 49999
$ws.Range("L56").Value = 49999
$ws.Range("N56").Value = -51689
$ws.Range("H58").Value = 2183.8948
$ws.Range("I58").Value = 1969.8
$ws.Range("J58").Value = 2421.7778
$ws.Range("K58").Value = 1969.8
$ws.Range("L58").Value = 2421.7778
$ws.Range("M58").Value = -1766.8
$ws.Range("N58").Value = -2827.7778
$ws.Range("H59").Value = 26783
$ws.Range("I59").Value = 7000
$ws.Range("J59").Value = 39971.668
$ws.Range("K59").Value = 7000
$ws.Range("L59").Value = 39971.668
$ws.Range("M59").Value = -5855
$ws.Range("N59").Value = -42261.668
$ws.Range("H60").Value = 30328.9
$ws.Range("I60").Value = 23833.334
$ws.Range("J60").Value = 40072.25
$ws.Range("K60").Value = 23833.334
$ws.Range("L60").Value = 40072.25
$ws.Range("M60").Value = -23322.334
$ws.Range("N60").Value = -41094.25
$ws.Range("H61").Value = 19099.092
$ws.Range("H62").Value = 7088
$ws.Range("I62").Value = 7284
$ws.Range("K62").Value = 7284
$ws.Range("M62").Value = -6660
$ws.Range("H65").Value = 7088
$ws.Range("I65").Value = 7284
$ws.Range("K65").Value = 36420
$ws.Range("M65").Value = -33300
$ws.Range("H68").Value = 67500
$ws.Range("J68").Value = 67500
$ws.Range("L68").Value = 67500
$ws.Range("N68").Value = -68998
$ws.Range("H71").Value = 67500
$ws.Range("J71").Value = 67500
$ws.Range("L71").Value = 202500
$ws.Range("N71").Value = -209988
$ws.Range("H74").Value = 60000
$ws.Range("J74").Value = 60000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -61748
$ws.Range("H77").Value = 60000
$ws.Range("J77").Value = 60000
$ws.Range("L77").Value = 180000
$ws.Range("N77").Value = -188736
$ws.Range("H99").Value = 3075.4
$ws.Range("J99").Value = 3008.25
$ws.Range("L99").Value = 3008.25
$ws.Range("N99").Value = -6004.25
$ws.Range("H126").Value = 3075.4
$ws.Range("J126").Value = 3008.25
$ws.Range("L126").Value = 9024.75
$ws.Range("N126").Value = -13964.75
$ws.Range("H134").Value = 3413.75
$ws.Range("J134").Value = 3405.2856
$ws.Range("L134").Value = 10215.8568
$ws.Range("N134").Value = -15285.8568
$ws.Range("H136").Value = 2183.8948
$ws.Range("I136").Value = 1969.8
$ws.Range("J136").Value = 2421.7778
$ws.Range("K136").Value = 5909.4
$ws.Range("L136").Value = 7265.3334
$ws.Range("M136").Value = -3359.4
$ws.Range("N136").Value = -12365.3334

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 1238441.5
$ws.Range("I4").Value = 1264935
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 3794805
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -3794693
$ws.Range("N4").Value = -3000224
$ws.Range("H13").Value = 6212.143
$ws.Range("I13").Value = 396.75
$ws.Range("J13").Value = 13966
$ws.Range("K13").Value = 1190.25
$ws.Range("L13").Value = 41898
$ws.Range("M13").Value = -1022.25
$ws.Range("N13").Value = -42234
$ws.Range("H100").Value = 22618.334
$ws.Range("I100").Value = 577.5
$ws.Range("K100").Value = 1732.5
$ws.Range("M100").Value = -921.5
$ws.Range("H131").Value = 2471063.8
$ws.Range("J131").Value = 3032462
$ws.Range("L131").Value = 9097386
$ws.Range("N131").Value = -9107466
$ws.Range("H137").Value = 4433.143
$ws.Range("I137").Value = 1891
$ws.Range("J137").Value = 5450
$ws.Range("K137").Value = 5673
$ws.Range("L137").Value = 16350
$ws.Range("M137").Value = -573
$ws.Range("N137").Value = -26550

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5446
$ws.Range("H80").Value = 12643.214
$ws.Range("J80").Value = 7666.6665
$ws.Range("L80").Value = 7666.6665
$ws.Range("N80").Value = -9662.666499999999
$ws.Range("H83").Value = 12643.214
$ws.Range("J83").Value = 7666.6665
$ws.Range("L83").Value = 38333.3325
$ws.Range("N83").Value = -48317.3325
$ws.Range("H105").Value = 49750
$ws.Range("J105").Value = 49750
$ws.Range("L105").Value = 49750
$ws.Range("N105").Value = -56738
$ws.Range("H132").Value = 3374.6667
$ws.Range("I132").Value = 3473.0715
$ws.Range("J132").Value = 1997
$ws.Range("K132").Value = 10419.2145
$ws.Range("L132").Value = 5991
$ws.Range("M132").Value = -7889.2145
$ws.Range("N132").Value = -11051

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1515.2
$ws.Range("I16").Value = 775.1818
$ws.Range("K16").Value = 775.1818
$ws.Range("M16").Value = -605.1818
$ws.Range("H68").Value = 11344.637
$ws.Range("I68").Value = 2250.5
$ws.Range("K68").Value = 2250.5
$ws.Range("M68").Value = -1501.5
$ws.Range("H71").Value = 11344.637
$ws.Range("I71").Value = 2250.5
$ws.Range("K71").Value = 11252.5
$ws.Range("M71").Value = -7508.5
$ws.Range("H82").Value = 5628.2856
$ws.Range("J82").Value = 8166.6665
$ws.Range("L82").Value = 8166.6665
$ws.Range("N82").Value = -8888.666499999999
$ws.Range("H85").Value = 5628.2856
$ws.Range("J85").Value = 8166.6665
$ws.Range("L85").Value = 8166.6665
$ws.Range("N85").Value = -10662.6665
$ws.Range("H122").Value = 5822
$ws.Range("I122").Value = 5365.5
$ws.Range("J122").Value = 6735
$ws.Range("K122").Value = 16096.5
$ws.Range("L122").Value = 20205
$ws.Range("M122").Value = -13646.5
$ws.Range("N122").Value = -25105
$ws.Range("H132").Value = 4229.185
$ws.Range("I132").Value = 4037.577
$ws.Range("K132").Value = 12112.731
$ws.Range("M132").Value = -9582.731

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H75").Value = 95950
$ws.Range("J75").Value = 95950
$ws.Range("L75").Value = 95950
$ws.Range("N75").Value = -97822
$ws.Range("H78").Value = 95950
$ws.Range("J78").Value = 95950
$ws.Range("L78").Value = 287850
$ws.Range("N78").Value = -297210
$ws.Range("H100").Value = 483.46667
$ws.Range("I100").Value = 307.3
$ws.Range("K100").Value = 614.6
$ws.Range("M100").Value = -73.60000000000002
$ws.Range("H122").Value = 5139.05
$ws.Range("I122").Value = 1973.125
$ws.Range("J122").Value = 7249.6665
$ws.Range("K122").Value = 5919.375
$ws.Range("L122").Value = 21748.9995
$ws.Range("M122").Value = -3469.375
$ws.Range("N122").Value = -26648.9995
$ws.Range("H124").Value = 75153.60000000001
$ws.Range("J124").Value = 88329.336
$ws.Range("L124").Value = 88329.336
$ws.Range("N124").Value = -98149.336
$ws.Range("H132").Value = 6160.839
$ws.Range("I132").Value = 4140.3184
$ws.Range("K132").Value = 12420.9552
$ws.Range("M132").Value = -9890.9552
$ws.Range("H136").Value = 2809.4546
$ws.Range("I136").Value = 2989.5715
$ws.Range("J136").Value = 2494.25
$ws.Range("K136").Value = 8968.7145
$ws.Range("L136").Value = 7482.75
$ws.Range("M136").Value = -6418.7145
$ws.Range("N136").Value = -12582.75
